# Apply the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.785.85"
$ws.Range("E2").Value = "  -1.98%  "

# Row 3
$ws.Range("D3").Value = "2.237.05"
$ws.Range("E3").Value = "  -2.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.60"
$ws.Range("E5").Value = "  -5.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "297.46"
$ws.Range("E6").Value = "  +11.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -2.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  -0.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.30"
$ws.Range("E10").Value = "  -6.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -1.62%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.81"
$ws.Range("E12").Value = "  +2.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.01"
$ws.Range("E13").Value = "  -2.42%  "

# Row 14
$ws.Range("E14").Value = "  -2.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.25"
$ws.Range("E15").Value = "  -1.93%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.887"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").Value = "2.574.88"
$ws.Range("E17").Value = "  -2.21%  "

# Row 18
$ws.Range("D18").Value = "2.263.82"
$ws.Range("E18").Value = "  -1.00%  "

# Row 19
$ws.Range("D19").Value = "42.638.48"
$ws.Range("E19").Value = "  -2.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  +7.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000107"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.10"
$ws.Range("E22").Value = "  +1.00%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.56"
$ws.Range("E23").Value = "  +23.43%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -6.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.00"
$ws.Range("E25").Value = "  -1.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  -3.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.97"
$ws.Range("E27").Value = "  +0.96%  "

# Row 28
$ws.Range("E28").Value = "  -1.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.53"
$ws.Range("E29").Value = "  -8.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -0.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.26"
$ws.Range("E31").Value = "  -3.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.18"
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.20"
$ws.Range("E33").Value = "  -2.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0892"
$ws.Range("E34").Value = "  -2.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  -1.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.04"
$ws.Range("E36").Value = "  +6.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  +10.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.127"
$ws.Range("E38").Value = "  -1.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0369"
$ws.Range("E39").Value = "  -2.91%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  -1.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.237"
$ws.Range("E42").Value = "  -0.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.89"
$ws.Range("E43").Value = "  -3.72%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("E44").Value = "  -8.21%  "

# Row 45
$ws.Range("E45").Value = "  +0.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.33"
$ws.Range("E46").Value = "  -3.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.58"
$ws.Range("E47").Value = "  -6.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.33"
$ws.Range("E48").Value = "  +3.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.35"
$ws.Range("E49").Value = "  +3.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0985"
$ws.Range("E51").Value = "  -1.83%  "

